$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2023-10-01 22:52:04"
$ws.Range("B3").Value = "aptos2019"
$ws.Range("C3").Value = 1.321817548342189
